$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update corrected values for existing rows 174-176 ---

# Row 174
$ws.Range("B174").Value = 7310
$ws.Range("C174").Value = 5888
$ws.Range("D174").Value = 16758
$ws.Range("E174").Value = 10870
$ws.Range("F174").Value = 1422
$ws.Range("G174").Value = 2386
$ws.Range("H174").Value = 964
$ws.Range("I174").Value = 30978
$ws.Range("J174").Value = 27629
$ws.Range("K174").Value = 154
$ws.Range("L174").Value = 1461
$ws.Range("M174").Value = 4766
$ws.Range("N174").Value = 2509
$ws.Range("O174").Value = 2756
$ws.Range("P174").Value = 1684
$ws.Range("Q174").Value = 14299
$ws.Range("R174").Value = 3349

# Row 175
$ws.Range("B175").Value = 6801
$ws.Range("C175").Value = 5407
$ws.Range("D175").Value = 16504
$ws.Range("E175").Value = 11096
$ws.Range("F175").Value = 1393
$ws.Range("G175").Value = 2353
$ws.Range("H175").Value = 960
$ws.Range("I175").Value = 30914
$ws.Range("J175").Value = 27600
$ws.Range("K175").Value = 130
$ws.Range("L175").Value = 1427
$ws.Range("M175").Value = 4819
$ws.Range("N175").Value = 2383
$ws.Range("O175").Value = 2708
$ws.Range("P175").Value = 1622
$ws.Range("Q175").Value = 14510
$ws.Range("R175").Value = 3314

# Row 176
$ws.Range("B176").Value = 5957
$ws.Range("C176").Value = 4545
$ws.Range("D176").Value = 15746
$ws.Range("E176").Value = 11202
$ws.Range("F176").Value = 1412
$ws.Range("G176").Value = 2353
$ws.Range("H176").Value = 941
$ws.Range("I176").Value = 30243
$ws.Range("J176").Value = 26948
$ws.Range("K176").Value = 151
$ws.Range("L176").Value = 1007
$ws.Range("M176").Value = 4441
$ws.Range("N176").Value = 2507
$ws.Range("O176").Value = 2685
$ws.Range("P176").Value = 1628
$ws.Range("Q176").Value = 14530
$ws.Range("R176").Value = 3295

# --- Append new row 177 with the new monthly data point ---

# Column A holds a "dd-mm-yyyy"-looking label that must stay plain text
# (matching the rest of the Serie column) instead of being auto-converted
# to a date serial number by the usual Value assignment. Route it through
# a text formula + paste-as-values so it lands as a shared string.
$ws.Range("A177").Formula = '="01-08-2021"'
$ws.Range("A177").Copy()
$ws.Range("A177").PasteSpecial(-4163)

$ws.Range("B177").Value = 5408
$ws.Range("C177").Value = 3989
$ws.Range("D177").Value = 15798
$ws.Range("E177").Value = 11809
$ws.Range("F177").Value = 1419
$ws.Range("G177").Value = 2353
$ws.Range("H177").Value = 935
$ws.Range("I177").Value = 30896
$ws.Range("J177").Value = 27608
$ws.Range("K177").Value = 163
$ws.Range("L177").Value = 1215
$ws.Range("M177").Value = 4334
$ws.Range("N177").Value = 2566
$ws.Range("O177").Value = 2741
$ws.Range("P177").Value = 1635
$ws.Range("Q177").Value = 14953
$ws.Range("R177").Value = 3288
